$p = $ppt.ActivePresentation
$th = $p.SlideMaster.Theme
Write-Host "=== Theme members ==="
Write-Host ($th | Get-Member | Format-List | Out-String)
